$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 45915
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

$ws.Range("B11").Value = "21,2065"
$ws.Range("C11").Value = "15,0136"
$ws.Range("D11").Value = "15,0136"
$ws.Range("E11").Value = "15,0136"
